$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.347.52"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.716.19"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'224.69"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'0.5299"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.06685"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").Value = "'0.2647"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").Value = "'20.91"
$ws.Range("E10").Value = "  -2.95%  "
$ws.Range("D11").Value = "'0.07703"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "'4.491"
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("D13").Value = "1.951.97"
$ws.Range("D14").Value = "1.714.20"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "'0.5798"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "0.0₅8193"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "'67.73"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "27.363.37"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").Value = "'220.37"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "'4.644"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("D22").Value = "'10.43"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("D23").Value = "'6.028"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "'145.52"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").Value = "'7.252"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("D29").Value = "'16.21"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").Value = "'0.05378"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("D33").Value = "'3.394"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("D35").Value = "'2.854"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").Value = "'0.9521"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("D38").Value = "'0.5889"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("D39").Value = "1.158.84"
$ws.Range("E39").Value = "  +10.74%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "'5.836"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("D42").Value = "'1.007"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "'0.8414"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "1.858.83"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("D47").Value = "'57.78"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "'0.4572"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("D49").Value = "'8.125"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  -1.06%  "
